# Adding round3 DBP runs and making current
# Process groups bottom-to-top so earlier row numbers stay valid as rows shift down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group: 2050 Plan (current row 113 in original layout) ---
$ws.Rows("113").Copy()
$ws.Rows("114").Insert()
$ws.Range("I113").ClearContents()
$ws.Range("U113").ClearContents()
$ws.Range("B114").Value = "2050_TM160_DBP_Plan_03"
$ws.Range("F114").Value = "Updated landuse/popsyn and small network updates"
$ws.Range("G114").Value = "M:\urban_modeling\baus\PBA50Plus\PBA50Plus_InclusionaryZoning_v2"
$ws.Range("H114").Value = "PBA50Plus_InclusionaryZoning_v2"
$ws.Range("J114").Value = "BlueprintNetworks_v17\net_2050_Blueprint"
$ws.Range("L114").Value = "https://app.asana.com/0/1204085012544660/1206926596430572/f"
$ws.Range("M114").Value = 18.42

# --- Group: 2050 No Project (current row 111 in original layout) ---
$ws.Rows("111").Copy()
$ws.Rows("112").Insert()
$ws.Range("I111").ClearContents()
$ws.Range("U111").ClearContents()
$ws.Range("B112").Value = "2050_TM160_DBP_NoProject_03"
$ws.Range("F112").Value = "Updated landuse/popsyn and small network updates"
$ws.Range("G112").Value = "M:\urban_modeling\baus\PBA50Plus\PBA50Plus_2020Validation_HE_BOC_v2"
$ws.Range("H112").Value = "PBA50Plus_2020Validation_HE_BOC_v2"
$ws.Range("J112").Value = "BlueprintNetworks_v17\net_2030_Baseline"
$ws.Range("L112").Value = "https://app.asana.com/0/1204085012544660/1206926596430574/f"
$ws.Range("M112").Value = 18.42

# --- Group: 2035 Plan (current row 108 in original layout) ---
$ws.Rows("108").Copy()
$ws.Rows("109").Insert()
$ws.Range("I108").ClearContents()
$ws.Range("U108").ClearContents()
$ws.Range("B109").Value = "2035_TM160_DBP_Plan_03"
$ws.Range("F109").Value = "Updated landuse/popsyn and small network updates"
$ws.Range("G109").Value = "M:\urban_modeling\baus\PBA50Plus\PBA50Plus_InclusionaryZoning_v2"
$ws.Range("H109").Value = "PBA50Plus_InclusionaryZoning_v2"
$ws.Range("J109").Value = "BlueprintNetworks_v17\net_2035_Blueprint"
$ws.Range("L109").Value = "https://app.asana.com/0/1204085012544660/1206926596430568/f"
$ws.Range("M109").Value = 16.82

# --- Group: 2035 No Project (current row 106 in original layout) ---
$ws.Rows("106").Copy()
$ws.Rows("107").Insert()
$ws.Range("I106").ClearContents()
$ws.Range("U106").ClearContents()
$ws.Range("B107").Value = "2035_TM160_DBP_NoProject_03"
$ws.Range("F107").Value = "Updated landuse/popsyn and small network updates"
$ws.Range("G107").Value = "M:\urban_modeling\baus\PBA50Plus\PBA50Plus_2020Validation_HE_BOC_v2"
$ws.Range("H107").Value = "PBA50Plus_2020Validation_HE_BOC_v2"
$ws.Range("J107").Value = "BlueprintNetworks_v17\net_2030_Baseline"
$ws.Range("L107").Value = "https://app.asana.com/0/1204085012544660/1206926596430570/f"
$ws.Range("M107").Value = 16.82

# --- Update sheet view: selection/frozen pane to reflect the new bottom area ---
$ws.Application.ActiveWindow.ScrollRow = 91
$ws.Range("K114").Select()
